# Weekly update: reorder the price rows (columns D, I, J, K, L, M, O, P)
# according to the mapping observed between the old and new snapshot.
# destRow -> sourceRow (i.e. the value that ends up in destRow comes from
# the data that used to live in sourceRow before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 27
    3  = 9
    4  = 10
    5  = 23
    6  = 24
    7  = 16
    8  = 21
    9  = 8
    10 = 35
    11 = 31
    12 = 32
    13 = 19
    14 = 20
    15 = 3
    16 = 11
    17 = 12
    18 = 28
    19 = 18
    20 = 29
    21 = 4
    22 = 5
    23 = 30
    24 = 34
    25 = 6
    26 = 7
    27 = 2
    28 = 36
    29 = 15
    30 = 14
    31 = 25
    32 = 26
    33 = 22
    34 = 13
    35 = 33
    36 = 17
}

$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# Snapshot current values for the affected columns across every data row
# (2..36) BEFORE writing anything, since source and destination rows
# overlap (this is a permutation of the existing rows).
$snapshot = @{}
for ($row = 2; $row -le 36; $row++) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $snapshot[$addr] = $ws.Range($addr).Value()
    }
}

# Now write the permuted values into their destination rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value = $snapshot[$srcAddr]
    }
}

$wb.Save()
